# Apply the "feat: add 2022-Q3 data" edit:
#  1. Insert a brand new worksheet named "2022-Q3" right after the "总计"
#     (summary) sheet, populated with the Q3-2022 fund holding data.
#  2. Update the "总计" summary sheet so its top data row now reports the
#     new 2022-Q3 figures, with every older quarter's row shifting down by
#     one and a new trailing row appearing for 2020-Q4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q3" worksheet right after "总计" (sheet 1)
# ---------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

$newSheet = $wb.Worksheets.Add($null, $summarySheet)
$newSheet.Name = "2022-Q3"

# Copy the header/body cell formatting from the existing 2022-Q2 sheet so
# the new sheet matches the look (borders, bold header, centered index
# column, etc.) of its siblings.
$q2Sheet.Range("A1:H5").Copy()
$newSheet.Range("A1:H5").PasteSpecial(-4122)  # xlPasteFormats

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "161017"
$newSheet.Range("C2").Value = "富国中证500指数增强（LOF）"
$newSheet.Range("D2").Value = "66.37"
$newSheet.Range("E2").Value = "90.18"
$newSheet.Range("F2").Value = "0.87"
$newSheet.Range("G2").Value = "0.5774"
$newSheet.Range("H2").Value = 4

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "014917"
$newSheet.Range("C3").Value = "汇丰晋信时代先锋混合A"
$newSheet.Range("D3").Value = "8.51"
$newSheet.Range("E3").Value = "94.48"
$newSheet.Range("F3").Value = "6.29"
$newSheet.Range("G3").Value = "0.5353"
$newSheet.Range("H3").Value = 3

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "014918"
$newSheet.Range("C4").Value = "汇丰晋信时代先锋混合C"
$newSheet.Range("D4").Value = "1.24"
$newSheet.Range("E4").Value = "94.48"
$newSheet.Range("F4").Value = "6.29"
$newSheet.Range("G4").Value = "0.0780"
$newSheet.Range("H4").Value = 3

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "013332"
$newSheet.Range("C5").Value = "富国中证500指数增强(LOF)C"
$newSheet.Range("D5").Value = "1.68"
$newSheet.Range("E5").Value = "90.18"
$newSheet.Range("F5").Value = "0.87"
$newSheet.Range("G5").Value = "0.0146"
$newSheet.Range("H5").Value = 4

# ---------------------------------------------------------------------
# Step 2: update the "总计" (summary) sheet with the new quarter on top
# ---------------------------------------------------------------------

# Extend the index column (A) down one more row, copying the existing
# formatting used by the rest of that column.
$summarySheet.Range("A8").Copy()
$summarySheet.Range("A9").PasteSpecial(-4122)  # xlPasteFormats
$summarySheet.Range("A9").Value = 7

$summarySheet.Range("B2").Value = "2022-Q3"
$summarySheet.Range("C2").Value = 4
$summarySheet.Range("D2").Value = 1.21

$summarySheet.Range("B3").Value = "2022-Q2"
$summarySheet.Range("C3").Value = 6
$summarySheet.Range("D3").Value = 2.94

$summarySheet.Range("B4").Value = "2022-Q1"
$summarySheet.Range("C4").Value = 5
$summarySheet.Range("D4").Value = 2.61

$summarySheet.Range("B5").Value = "2021-Q4"
$summarySheet.Range("C5").Value = 5
$summarySheet.Range("D5").Value = 1.83

$summarySheet.Range("B6").Value = "2021-Q3"
$summarySheet.Range("C6").Value = 24
$summarySheet.Range("D6").Value = 9.53

$summarySheet.Range("B7").Value = "2021-Q2"
$summarySheet.Range("C7").Value = 11
$summarySheet.Range("D7").Value = 1.15

$summarySheet.Range("B8").Value = "2021-Q1"
$summarySheet.Range("C8").Value = 13
$summarySheet.Range("D8").Value = 0.75

$summarySheet.Range("B9").Value = "2020-Q4"
$summarySheet.Range("C9").Value = 3
$summarySheet.Range("D9").Value = 0.11

Write-Host "2022-Q3 sheet added and 总计 sheet updated"
